$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "exitoso"
$ws.Range("B3").Value = "exitoso"
$ws.Range("B4").Value = "exitoso"
